# code updte cfor comps
# Add a third experiment-model comparison sheet (EM_Exp3) to the workbook,
# re-using the two "Est_Prop*" configuration rows that already exist on
# EM_Exp2 (Est_PropWith_SR_Y / Est_PropAcr), and leave a couple of blank
# rows below them for future additions.

$wb = $excel.ActiveWorkbook

# --- EM_Exp2: cursor was left on row 7 and "Select All" was pressed, so the
# highlighted range grows to the whole sheet while the active cell stays
# put on A7 ------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Select() | Out-Null

# --- add the new sheet after EM_Exp2 -----------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "EM_Exp3"

# column A is a bit narrower than on EM_Exp2
$ws3.Columns.Item(1).ColumnWidth = 32.17

# give the whole block a plain black Calibri font (new style vs. the
# Times New Roman header style used on the other sheets)
$ws3.Range("A1:K5").Font.Color = 0

# header row, identical layout/order to EM_Exp2
$headers = @("EM_Name","use_fish_sexRatio","use_srv_sexRatio","fish_age_prop","srv_age_prop","fish_len_prop","srv_len_prop","est_sexRatio_par","share_M_sex","sexRatio_al_or_y","SR_Wt")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws3.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# row 2: Est_PropWith_SR_Y
$ws3.Range("A2").Value = "Est_PropWith_SR_Y"
$ws3.Range("B2").Value = $true
$ws3.Range("C2").Value = $true
$ws3.Range("D2").Value = "within"
$ws3.Range("E2").Value = "within"
$ws3.Range("F2").Value = "within"
$ws3.Range("G2").Value = "within"
$ws3.Range("H2").Value = $true
$ws3.Range("I2").Value = $false
$ws3.Range("J2").Value = "within_year_only"
$ws3.Range("K2").Value = 1

# row 3: Est_PropAcr
$ws3.Range("A3").Value = "Est_PropAcr"
$ws3.Range("B3").Value = $false
$ws3.Range("C3").Value = $false
$ws3.Range("D3").Value = "across"
$ws3.Range("E3").Value = "across"
$ws3.Range("F3").Value = "across"
$ws3.Range("G3").Value = "across"
$ws3.Range("H3").Value = $true
$ws3.Range("I3").Value = $false
$ws3.Range("J3").Value = "None"
$ws3.Range("K3").Value = 0

# rows 4-5 stay blank, reserved for more rows later; leave the cursor
# selecting them, which is where the sheet was left/saved
$ws3.Range("A4:A5").EntireRow.Select() | Out-Null
